$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New constant values for columns B (time in ms) and C (auto capacity) across all data rows (2-55)
$newB = 33.94444444444444
$newC = 1.95

# Updated D (auto scs / dic_nbre_clients_poisson_2_keys) and E (probability) values for existing rows 2-51
$dVals = @(0,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,49,51,52)
$eVals = @(0.126,0.002,0.005,0.013,0.026,0.039,0.044,0.051,0.05,0.029,0.032,0.027,0.026,0.029,0.03,0.033,0.037,0.045,0.032,0.026,0.023,0.027,0.02,0.032,0.011,0.022,0.024,0.017,0.008,0.02,0.015,0.01,0.007,0.01,0.01,0.006,0.002,0.001,0.002,0.003,0.003,0.004,0.003,0.003,0.002,0.001,0.002,0.001,0.002,0.002)

for ($i = 0; $i -lt $dVals.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $newB
    $ws.Cells.Item($r, 3).Value = $newC
    $ws.Cells.Item($r, 4).Value = $dVals[$i]
    $ws.Cells.Item($r, 5).Value = $eVals[$i]
}

# Append the new rows 52-55 (A, D, E vary; B, C constant as above)
# Column A on these rows reuses the existing bold/centered/bordered style from A51 (style index referenced by A2:A51)
$styleSource = $ws.Range("A51")

$newA = @(50,51,52,53)
$newD = @(55,57,61,63)
$newE = @(0.001,0.001,0.001,0.001)

for ($i = 0; $i -lt $newA.Length; $i++) {
    $r = 52 + $i
    $aCell = $ws.Cells.Item($r, 1)
    $styleSource.Copy()
    $aCell.PasteSpecial(-4122)
    $aCell.Value = $newA[$i]
    $ws.Cells.Item($r, 2).Value = $newB
    $ws.Cells.Item($r, 3).Value = $newC
    $ws.Cells.Item($r, 4).Value = $newD[$i]
    $ws.Cells.Item($r, 5).Value = $newE[$i]
}

$excel.CutCopyMode = $false

# Dimension should now span A1:E55 (Excel recalculates this automatically on save, but we confirm the used range)
